$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 868; this shifts every
# existing row from 868..909 down by one (869..910), which matches the
# bulk of the diff (each old row's data reappears one row lower).
$ws.Rows(868).Insert()

# Populate the newly inserted row 868 with its data. Columns A/B hold
# text (date-as-text / weekday-as-text) so force a text number format
# before assigning, then drop back to the Normal style so no stray
# style index is left on the cell (matches surrounding data cells that
# carry no explicit style).
$ws.Range("A868").NumberFormat = "@"
$ws.Range("A868").Value = "2026/02/23"
$ws.Range("A868").Style = "Normal"

$ws.Range("B868").NumberFormat = "@"
$ws.Range("B868").Value = "月"
$ws.Range("B868").Style = "Normal"

$ws.Range("C868").Value = 19
$ws.Range("D868").Value = 201
